# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 15 (pushing the existing
# rows 15-79 down to 16-80) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..79 down by one row, leaving a new blank row 15.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new data point.
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44624
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 100112021
$ws.Range("G15").Value = "Ají"
$ws.Range("H15").Value = "Chilena(o)"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 220
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 23000
$ws.Range("M15").Value = 22455
$ws.Range("N15").Value = "$/caja 25 kilos"
$ws.Range("O15").Value = "Provincia de Huasco"
$ws.Range("P15").Value = 898
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
